$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Years" header in row 1 by deleting the whole row; this shifts
# the year values (1950..2008) up by one row so the sheet goes from
# A1:A60 (Years, 1950..2008) to A1:A59 (1950..2008), and also drops the
# now-unused "Years" shared string.
$ws.Rows.Item(1).Delete()

# Update the selected cell to match the target workbook's saved selection.
$ws.Range("D9").Select()
